# Dialogs.xlsx edit: "Almost done with Rename (Type, Method, Event). SS"
#
# Summary of changes applied:
#  1. Sheet2 becomes the active/selected tab (was Sheet1).
#  2. Sheet2's view selection moves to P11 (was C34).
#  3. Four new cell comments are added to Sheet2 (S10, S12, S14, S16) by
#     Gerald Rubin documenting the status of the Rename Type/Method/Event work.
#  4. A new status-tracking grid is added to Sheet2 in columns K:S, rows 9-20,
#     capturing progress on Add/Rename Type, Method, Property and Event work.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1 & 2: activate Sheet2 (updates workbookView.activeTab + tabSelected on the
# sheets) and move the on-screen selection to P11.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("P11").Select()

# ---------------------------------------------------------------------------
# 4: New status grid, columns K (11) .. S (19)
# ---------------------------------------------------------------------------

$xlCenter = -4108

# --- Row 9: wrapped, centered column headers (row height taller to fit) ---
$ws2.Rows.Item(9).RowHeight = 43.5
$headerCols = @(12,13,14,15,16,17,18,19)
$headerText = @(
    "Call Client check for duplicate in Save func",
    "Client routine written",
    "Types routine written",
    "Comics routine written",
    "Replace in array",
    "Replace in code",
    "Add to array",
    "Add to code"
)
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $ws2.Cells.Item(9, $headerCols[$i])
    $cell.Value = $headerText[$i]
    $cell.HorizontalAlignment = $xlCenter
    $cell.WrapText = $true
}

# --- Column K (11): row labels for each action being tracked ---
$kLabels = @{
    10 = "Add Type";
    11 = "Rename Type";
    12 = "Add Method";
    13 = "Rename Method";
    14 = "Add Property";
    15 = "Edit Property";
    16 = "Add Event";
    17 = "Rename Event";
}
foreach ($row in $kLabels.Keys) {
    $ws2.Cells.Item($row, 11).Value = $kLabels[$row]
}

# --- Columns L:S (12-19), rows 10-17: the "x" / "-" completion grid ---
# Each entry: row -> column -> (value, wrapText?)
# value "x" = done, "-" = not applicable / not done, "incomplete" = WIP.
$grid = @{
    10 = @{ 12="x"; 13="x"; 14="-"; 15="x"; 16="-"; 17="-"; 18="x"; 19="x" };
    11 = @{ 12="x"; 13="x"; 14="-"; 15="x";                18="-"; 19="-" };
    12 = @{ 12="x"; 13="x"; 14="x"; 15="-"; 16="-"; 17="-"; 18="x"; 19="x" };
    13 = @{ 12="x"; 13="x"; 14="x"; 15="-";                18="-"; 19="-" };
    14 = @{ 12="x"; 13="x"; 14="x"; 15="-"; 16="-"; 17="-"; 18="x"; 19="x" };
    15 = @{ 12="x"; 13="x"; 14="x"; 15="-";                18="-"; 19="-" };
    16 = @{ 12="x"; 13="x"; 14="x"; 15="-"; 16="-"; 17="-"; 18="x"; 19="incomplete" };
    17 = @{ 12="x"; 13="x"; 14="x"; 15="-";                18="-"; 19="-" };
}
# Columns that carry the softer "quote-prefixed" look in the original
# workbook (M, N always wrap; L, O, P, Q, R, S stay single-line except M/N).
$wrapCols = @(13, 14)

foreach ($row in ($grid.Keys | Sort-Object)) {
    $rowData = $grid[$row]
    foreach ($col in ($rowData.Keys | Sort-Object)) {
        $val = $rowData[$col]
        $cell = $ws2.Cells.Item($row, $col)
        $cell.Value = $val
        $cell.HorizontalAlignment = $xlCenter
        if ($wrapCols -contains [int]$col) {
            $cell.WrapText = $true
        }
    }
}

# S16 ("incomplete") carries a stray date-ish number format (numFmtId 16,
# "d-mmm") in the source workbook, left over from a formatting accident.
$ws2.Range("S16").NumberFormat = "d-mmm"

# --- Row 19/20: a couple of standalone "Is it correct?" / "x" cells ---
$ws2.Cells.Item(19, 14).Value = "Is it correct?"
$ws2.Cells.Item(19, 14).WrapText = $true
$ws2.Cells.Item(19, 15).Value = "Is it correct?"
$ws2.Cells.Item(19, 15).HorizontalAlignment = $xlCenter

$ws2.Cells.Item(20, 14).Value = "x"
$ws2.Cells.Item(20, 14).WrapText = $true
$ws2.Cells.Item(20, 15).Value = "x"
$ws2.Cells.Item(20, 15).HorizontalAlignment = $xlCenter

# --- Column widths for the new columns ---
$ws2.Columns.Item(11).ColumnWidth = 14.81640625
$ws2.Columns.Item(12).ColumnWidth = 16.81640625
$ws2.Columns.Item(13).ColumnWidth = 13.453125
$ws2.Columns.Item(14).ColumnWidth = 14.54296875
$ws2.Columns.Item(15).ColumnWidth = 13.36328125
$ws2.Columns.Item(16).ColumnWidth = 8.7265625
$ws2.Columns.Item(17).ColumnWidth = 8.7265625
$ws2.Columns.Item(18).ColumnWidth = 8.7265625
$ws2.Columns.Item(19).ColumnWidth = 10.453125

# ---------------------------------------------------------------------------
# 3: Comments documenting where each Rename/Add routine lives in the code.
# ---------------------------------------------------------------------------
$c1 = $ws2.Range("S10").AddComment("Gerald Rubin:" + [char]10 + "In Types.js#addItem")
$c1.Author = "Gerald Rubin"

$c2 = $ws2.Range("S12").AddComment("Gerald Rubin:" + [char]10 + "In Client.js#addMethodToActiveType")
$c2.Author = "Gerald Rubin"

$c3 = $ws2.Range("S14").AddComment("Gerald Rubin:" + [char]10 + "But Code.js#m_functionAdd_Type_Property needs examination by Ken")
$c3.Author = "Gerald Rubin"

$c4 = $ws2.Range("S16").AddComment("Gerald Rubin:" + [char]10 + "In Client.js#addEventToActiveType, but routine is commented out in Code.js#m_functionAdd_Type_Event")
$c4.Author = "Gerald Rubin"

Write-Output "Dialogs.xlsx update applied."
